# Apply the "new TPM" data refresh:
#  - remove the "ECs" sending-cluster row (old row 2)
#  - update the remaining two rows (FAPs, MuSCs) with recomputed values

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the whole "ECs" row (row 2); rows below shift up, shared strings /
# dimension are re-written automatically on save.
$ws.Rows.Item(2).Delete()

# After the delete, the former row 3 ("FAPs") is now row 2, and the former
# row 4 ("MuSCs") is now row 3. Their label columns (A:D) already contain
# the correct values, only the numeric columns (E:T) need to be refreshed.

# New row 2 ("FAPs" -> Dsc2 -> Dsg1a -> MuSCs)
$ws.Cells.Item(2, 5).Value  = 1                        # E2
$ws.Cells.Item(2, 6).Value  = 0.3333333333333333        # F2
$ws.Cells.Item(2, 7).Value  = 0.008366333333333333      # G2
$ws.Cells.Item(2, 8).Value  = 0.025099                  # H2
$ws.Cells.Item(2, 9).Value  = 0.1196620722863995        # I2
$ws.Cells.Item(2, 10).Value = 0.1693601171397918        # J2
$ws.Cells.Item(2, 11).Value = 1                         # K2
$ws.Cells.Item(2, 12).Value = 0.5                       # L2
$ws.Cells.Item(2, 13).Value = 0.007882                  # M2
$ws.Cells.Item(2, 14).Value = 0.015764                  # N2
$ws.Cells.Item(2, 15).Value = 1                         # O2
$ws.Cells.Item(2, 16).Value = 1                         # P2
$ws.Cells.Item(2, 17).Value = 0.00006594343933333333    # Q2
$ws.Cells.Item(2, 18).Value = 0.000395660636            # R2
$ws.Cells.Item(2, 19).Value = 0.1196620722863995        # S2
$ws.Cells.Item(2, 20).Value = 0.1693601171397918        # T2

# New row 3 ("MuSCs" -> Dsc2 -> Dsg1a -> MuSCs)
$ws.Cells.Item(3, 5).Value  = 1                         # E3
$ws.Cells.Item(3, 6).Value  = 0.5                       # F3
$ws.Cells.Item(3, 7).Value  = 0.06155                   # G3
$ws.Cells.Item(3, 8).Value  = 0.1231                    # H3
$ws.Cells.Item(3, 9).Value  = 0.8803379277136005        # I3
$ws.Cells.Item(3, 10).Value = 0.8306398828602083        # J3
$ws.Cells.Item(3, 11).Value = 1                         # K3
$ws.Cells.Item(3, 12).Value = 0.5                       # L3
$ws.Cells.Item(3, 13).Value = 0.007882                  # M3
$ws.Cells.Item(3, 14).Value = 0.015764                  # N3
$ws.Cells.Item(3, 15).Value = 1                         # O3
$ws.Cells.Item(3, 16).Value = 1                         # P3
$ws.Cells.Item(3, 17).Value = 0.0004851371              # Q3
$ws.Cells.Item(3, 18).Value = 0.0019405484              # R3
$ws.Cells.Item(3, 19).Value = 0.8803379277136005        # S3
$ws.Cells.Item(3, 20).Value = 0.8306398828602083        # T3

$wb.Save()
